$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 358.5
$ws.Range("I38").Value = 358.5
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 1075.5
$ws.Range("L38").Value = 0
$ws.Range("N38").Value = -703.5
$ws.Range("M38").ClearContents()

$ws.Range("H52").Value = 3914.9
$ws.Range("I52").Value = 3329
$ws.Range("J52").Value = 3980
$ws.Range("K52").Value = 9987
$ws.Range("L52").Value = 11940
$ws.Range("M52").Value = -9827
$ws.Range("N52").Value = -12260

$ws.Range("H62").Value = 2399.6667
$ws.Range("I62").Value = 2099.5
$ws.Range("K62").Value = 2099.5
$ws.Range("M62").Value = -1475.5

$ws.Range("H65").Value = 2399.6667
$ws.Range("I65").Value = 2099.5
$ws.Range("K65").Value = 10497.5
$ws.Range("M65").Value = -7377.5

$ws.Range("H125").Value = 1195.5555
$ws.Range("J125").Value = 1065
$ws.Range("L125").Value = 9585
$ws.Range("N125").Value = -14505

$ws.Range("H132").Value = 1287
$ws.Range("I132").Value = 1026.0416
$ws.Range("J132").Value = 2852.75
$ws.Range("K132").Value = 3078.1248
$ws.Range("L132").Value = 8558.25
$ws.Range("M132").Value = -548.1248000000001
$ws.Range("N132").Value = -13618.25

$ws.Range("H137").Value = 1599.4
$ws.Range("I137").Value = 1599.4
$ws.Range("K137").Value = 4798.200000000001
$ws.Range("M137").Value = -2248.200000000001

$ws.Range("H138").Value = 2305.068
$ws.Range("J138").Value = 2578.3076
$ws.Range("L138").Value = 7734.9228
$ws.Range("N138").Value = -18014.9228

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6859.885
$ws.Range("I32").Value = 4441.6875
$ws.Range("J32").Value = 10729
$ws.Range("K32").Value = 4441.6875
$ws.Range("L32").Value = 10729
$ws.Range("M32").Value = -4154.6875
$ws.Range("N32").Value = -11303

$ws.Range("H45").Value = 1552
$ws.Range("I45").Value = 996.25
$ws.Range("K45").Value = 996.25
$ws.Range("M45").Value = -619.25

$ws.Range("H61").Value = 4303.0835
$ws.Range("I61").Value = 2401.8333
$ws.Range("K61").Value = 2401.8333
$ws.Range("M61").Value = -2189.8333

$ws.Range("H74").Value = 1374.9565
$ws.Range("I74").Value = 1241.5
$ws.Range("K74").Value = 1241.5
$ws.Range("M74").Value = -367.5

$ws.Range("H77").Value = 1374.9565
$ws.Range("I77").Value = 1241.5
$ws.Range("K77").Value = 6207.5
$ws.Range("M77").Value = -1839.5

$ws.Range("H96").Value = 40000
$ws.Range("J96").Value = 40000
$ws.Range("L96").Value = 40000
$ws.Range("N96").Value = -45492

$ws.Range("H105").Value = 49999.5
$ws.Range("J105").Value = 49999.5
$ws.Range("L105").Value = 49999.5
$ws.Range("N105").Value = -56987.5

$ws.Range("H122").Value = 650
$ws.Range("I122").Value = 650
$ws.Range("K122").Value = 1950
$ws.Range("M122").Value = 500

$ws.Range("H132").Value = 1816.5128
$ws.Range("I132").Value = 1258.6522
$ws.Range("K132").Value = 3775.9566
$ws.Range("M132").Value = -1245.9566

$ws.Range("H135").Value = 19294.334
$ws.Range("J135").Value = 19294.334
$ws.Range("L135").Value = 19294.334
$ws.Range("N135").Value = -29434.334

$ws.Range("H136").Value = 4303.0835
$ws.Range("I136").Value = 2401.8333
$ws.Range("K136").Value = 7205.499899999999
$ws.Range("M136").Value = -4655.499899999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H75").Value = 11999.667
$ws.Range("I75").Value = 11999.5
$ws.Range("K75").Value = 11999.5
$ws.Range("M75").Value = -11063.5

$ws.Range("H78").Value = 11999.667
$ws.Range("I78").Value = 11999.5
$ws.Range("K78").Value = 35998.5
$ws.Range("M78").Value = -31318.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2053.516
$ws.Range("I31").Value = 1586.8235
$ws.Range("K31").Value = 1586.8235
$ws.Range("M31").Value = -1291.8235

$ws.Range("H34").Value = 2053.516
$ws.Range("I34").Value = 1586.8235
$ws.Range("K34").Value = 1586.8235
$ws.Range("M34").Value = -1384.8235

$ws.Range("H74").Value = 30379.166
$ws.Range("J74").Value = 30379.166
$ws.Range("L74").Value = 30379.166
$ws.Range("N74").Value = -32127.166

$ws.Range("H77").Value = 30379.166
$ws.Range("J77").Value = 30379.166
$ws.Range("L77").Value = 91137.49800000001
$ws.Range("N77").Value = -99873.49800000001

$ws.Range("H134").Value = 1860.5667
$ws.Range("I134").Value = 1628.4445
$ws.Range("K134").Value = 4885.333500000001
$ws.Range("M134").Value = -2350.333500000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 6117338.5
$ws.Range("I4").Value = 162
$ws.Range("K4").Value = 486
$ws.Range("M4").Value = -374

$ws.Range("H134").Value = 2309
$ws.Range("I134").Value = 1746.3334
$ws.Range("J134").Value = 3997
$ws.Range("K134").Value = 5239.0002
$ws.Range("L134").Value = 11991
$ws.Range("M134").Value = -169.0002000000004
$ws.Range("N134").Value = -22131

$ws.Range("H139").Value = 10981.363
$ws.Range("I139").Value = 13986.875
$ws.Range("K139").Value = 41960.625
$ws.Range("M139").Value = -36820.625

$ws.Range("H140").Value = 1919.6061
$ws.Range("I140").Value = 1136.4667
$ws.Range("J140").Value = 2572.2222
$ws.Range("K140").Value = 3409.4001
$ws.Range("L140").Value = 7716.6666
$ws.Range("M140").Value = 1770.5999
$ws.Range("N140").Value = -18076.6666

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2899
$ws.Range("J80").Value = 2900
$ws.Range("L80").Value = 2900
$ws.Range("N80").Value = -4896

$ws.Range("H83").Value = 2899
$ws.Range("J83").Value = 2900
$ws.Range("L83").Value = 14500
$ws.Range("N83").Value = -24484

$ws.Range("H113").Value = 955.1818
$ws.Range("I113").Value = 684.8570999999999
$ws.Range("J113").Value = 1428.25
$ws.Range("K113").Value = 684.8570999999999
$ws.Range("L113").Value = 1428.25
$ws.Range("M113").Value = 1485.1429
$ws.Range("N113").Value = -5768.25

$ws.Range("H122").Value = 1466.8422
$ws.Range("I122").Value = 1144.8
$ws.Range("K122").Value = 3434.4
$ws.Range("M122").Value = -984.3999999999996

$ws.Range("H132").Value = 4278448
$ws.Range("I132").Value = 6414339
$ws.Range("K132").Value = 19243017
$ws.Range("M132").Value = -19240487

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3273.2632
$ws.Range("I7").Value = 2093.3076
$ws.Range("J7").Value = 5829.8335
$ws.Range("K7").Value = 2093.3076
$ws.Range("L7").Value = 5829.8335
$ws.Range("M7").Value = -1981.3076
$ws.Range("N7").Value = -6053.8335

$ws.Range("H126").Value = 3273.2632
$ws.Range("I126").Value = 2093.3076
$ws.Range("J126").Value = 5829.8335
$ws.Range("K126").Value = 6279.9228
$ws.Range("L126").Value = 17489.5005
$ws.Range("M126").Value = -3809.9228
$ws.Range("N126").Value = -22429.5005

$ws.Range("H132").Value = 3274.5
$ws.Range("I132").Value = 2586.8333
$ws.Range("J132").Value = 4099.7
$ws.Range("K132").Value = 7760.499899999999
$ws.Range("L132").Value = 12299.1
$ws.Range("M132").Value = -5230.499899999999
$ws.Range("N132").Value = -17359.1

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 1999
$ws.Range("I62").Value = 1999
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 1999
$ws.Range("L62").Value = 0
$ws.Range("N62").Value = -1375
$ws.Range("M62").ClearContents()

$ws.Range("H65").Value = 1999
$ws.Range("I65").Value = 1999
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 9995
$ws.Range("L65").Value = 0
$ws.Range("N65").Value = -6875
$ws.Range("M65").ClearContents()

$ws.Range("H75").Value = 30000
$ws.Range("J75").Value = 30000
$ws.Range("L75").Value = 30000
$ws.Range("N75").Value = -31872

$ws.Range("H78").Value = 30000
$ws.Range("J78").Value = 30000
$ws.Range("L78").Value = 90000
$ws.Range("N78").Value = -99360

$ws.Range("H107").Value = 1081.75
$ws.Range("I107").Value = 933.6667
$ws.Range("J107").Value = 1170.6
$ws.Range("K107").Value = 2801.0001
$ws.Range("L107").Value = 3511.8
$ws.Range("M107").Value = -881.0001000000002
$ws.Range("N107").Value = -7351.799999999999

$ws.Range("H126").Value = 2574.5715
$ws.Range("I126").Value = 1442.0667
$ws.Range("J126").Value = 5405.8335
$ws.Range("K126").Value = 4326.2001
$ws.Range("L126").Value = 16217.5005
$ws.Range("M126").Value = -1856.2001
$ws.Range("N126").Value = -21157.5005

$ws.Range("H132").Value = 2042.8572
$ws.Range("I132").Value = 1114.5714
$ws.Range("J132").Value = 2971.1428
$ws.Range("K132").Value = 3343.7142
$ws.Range("L132").Value = 8913.428400000001
$ws.Range("M132").Value = -813.7142000000003
$ws.Range("N132").Value = -13973.4284

$ws.Range("H135").Value = 121630
$ws.Range("J135").Value = 121630
$ws.Range("L135").Value = 121630
$ws.Range("N135").Value = -131770

$ws.Range("H136").Value = 16341802
$ws.Range("I136").Value = 23149570
$ws.Range("J136").Value = 3158.8
$ws.Range("K136").Value = 69448710
$ws.Range("L136").Value = 9476.400000000001
$ws.Range("M136").Value = -69446160
$ws.Range("N136").Value = -14576.4
